# Generate Report for Handback
# Populates the "297a6e8e-b874-44ea-b854-0074a6886a84" handback row (row 6)
# on both the zh-cn and de-de sheets with the freshly generated target file,
# handback file, handback datetime and an "out of date" error detail message,
# plus widens the Error Detail column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23a9dfc715cc7dc7c1e3a7ca061539b9c3c31f31/e2e/297a6e8e-b874-44ea-b854-0074a6886a84.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b746cd351562c1606fc13689f53cb74980b82d74/e2e/297a6e8e-b874-44ea-b854-0074a6886a84.md."

function Update-HandbackRow {
    param(
        [string]$SheetName,
        [string]$TargetFileValue,
        [string]$HandbackFileValue,
        [string]$HandbackDateTimeValue,
        [string]$HyperlinkTarget
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P) so the new message is readable.
    $ws.Columns.Item(16).ColumnWidth = 40

    # Latest Target File (I6) - styled + linked like column A's handback link.
    $ws.Range("I6").Value = $TargetFileValue
    $ws.Hyperlinks.Add($ws.Range("I6"), $HyperlinkTarget, "", "", $TargetFileValue) | Out-Null
    $ws.Range("I6").Style = $ws.Range("A6").Style

    # Latest Handback File (J6)
    $ws.Range("J6").Value = $HandbackFileValue

    # Latest Handback DateTime (K6)
    $ws.Range("K6").Value = $HandbackDateTimeValue

    # Error Detail (P6)
    $ws.Range("P6").Value = $errorDetail
}

Update-HandbackRow "zh-cn" `
    "297a6e8e-b874-44ea-b854-0074a6886a84.md" `
    "297a6e8e-b874-44ea-b854-0074a6886a84.e98931a284fd06bb7f2ecc41ff1425cbcd47c368.zh-cn.xlf" `
    "2016-11-29 03:03:12" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b746cd351562c1606fc13689f53cb74980b82d74/e2e/297a6e8e-b874-44ea-b854-0074a6886a84.md"

Update-HandbackRow "de-de" `
    "297a6e8e-b874-44ea-b854-0074a6886a84.md" `
    "297a6e8e-b874-44ea-b854-0074a6886a84.e98931a284fd06bb7f2ecc41ff1425cbcd47c368.de-de.xlf" `
    "2016-11-29 03:03:30" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b746cd351562c1606fc13689f53cb74980b82d74/e2e/297a6e8e-b874-44ea-b854-0074a6886a84.md"
